$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "p_sup [kN/m3]"
$ws.Range("D1").Value = "p_inf [kN/m3]"
